# Generate Report for Handback
# Refresh the handoff/handback correspondence timestamps for each locale
# sheet (the report is regenerated each time a handback is processed, so
# the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns advance to the latest run's timestamps).

$wb = $excel.ActiveWorkbook

# zh-cn sheet: both data rows share the same handoff/handback file, so both
# rows receive the same refreshed timestamps.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-18 08:33:17"
$wsZhCn.Range("D3").Value = "2016-02-18 08:33:17"
$wsZhCn.Range("G2").Value = "2016-02-18 08:34:12"
$wsZhCn.Range("G3").Value = "2016-02-18 08:34:12"

# de-de sheet: same refresh for its handoff/handback timestamps.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-18 08:33:29"
$wsDeDe.Range("D3").Value = "2016-02-18 08:33:29"
$wsDeDe.Range("G2").Value = "2016-02-18 08:34:43"
$wsDeDe.Range("G3").Value = "2016-02-18 08:34:43"
